# The commit swaps the presentation's theme colours: the slide master's
# theme (ppt/theme/theme1.xml, originally the "Integral" palette) ends up
# holding the stock "Office Theme" palette (the colours that used to live
# only in the notes-master's theme, ppt/theme/theme2.xml).
#
# PowerPoint's object model exposes the live (slide-master-facing) theme
# colours through Theme.ThemeColorScheme.Colors(i).RGB, in the standard
# clrScheme order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# so we recolour the theme in place to match the target "Office" palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$tcs = $theme.ThemeColorScheme

function Get-RGBValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = stock Office theme colours (clrScheme name="Office").
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = Get-RGBValue $officeThemeColors[$i - 1]
}
